# Update the workbook to match the target revision:
#  1. Column C (Förändrad) for all existing data rows (2..308) moves from
#     2023-09-21 (45190) to 2023-09-23 (45192).
#  2. Row 308 picks up an explicit row height (it is no longer the last row).
#  3. Three new announcement rows are appended (309, 310, 311).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bulk-update the "Förändrad" column for the existing rows.
$ws.Range("C2:C308").Value = 45192

# 2) Row 308 is no longer the last row in the sheet, so it gets the same
#    explicit "15pt custom height" stamp as every other data row.
$ws.Rows.Item(308).RowHeight = 15

# 3) Append the three new rows reported by the source feed.
$newRows = @(
    @{ Row = 309; A = "A 44819-2023"; B = 45190; C = 45192; G = 0.5 },
    @{ Row = 310; A = "A 45044-2023"; B = 45191; C = 45192; G = 0.5 },
    @{ Row = 311; A = "A 45085-2023"; B = 45191; C = 45192; G = 0.9 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "ÖSTERGÖTLANDS LÄN"
    $ws.Cells.Item($row, 5).Value = "SÖDERKÖPING"

    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0

    $ws.Cells.Item($row, 18).Value = ""
    $ws.Cells.Item($row, 18).WrapText = $true
}

# Rows 309 and 310 also get the explicit "15pt custom height" stamp, just
# like every other data row except the new final row (311).
$ws.Rows.Item(309).RowHeight = 15
$ws.Rows.Item(310).RowHeight = 15
